$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version 3.8.0 -> 3.9.0
$ws.Cells.Item(3, 2).Value = "3.9.0"

# Experimental now has an explicit value of "false".
# A plain Value assignment of the literal text "false" would be auto-coerced
# to a boolean cell by the engine, so build it via a helper cell holding a
# text formula, then copy/paste only the values into the target cell. This
# keeps the cell as a shared-string text cell (matching the source data)
# and does not disturb the existing cell style.
$helper = $ws.Cells.Item(100, 20)
$helper.Formula = "=""false"""
$helper.Copy()
$ws.Cells.Item(7, 2).PasteSpecial(-4163)
$helper.Clear()

# Date updated
$ws.Cells.Item(8, 2).Value = "2024-12-02T17:05:26-06:00"

# Contact rows updated (3 rows now carry distinct contact detail strings)
$ws.Cells.Item(10, 2).Value = "null (https://www.ihe.net/ihe_domains/it_infrastructure/)"
$ws.Cells.Item(11, 2).Value = "null (iti@ihe.net)"
$ws.Cells.Item(12, 2).Value = "IHE IT Infrastructure Technical Committee (iti@ihe.net)"

# Jurisdiction updated
$ws.Cells.Item(13, 2).Value = "Global (Whole world)"

$excel.CutCopyMode = 0
